$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new worksheet "IAM029" right before "Test Case Steps"
#    (matches sheetId=14 / rId11, with "Test Case Steps" shifting to rId12)
# ------------------------------------------------------------------
$testCaseSteps = $wb.Worksheets.Item("Test Case Steps")
$iam029 = $wb.Worksheets.Add($testCaseSteps)
$iam029.Name = "IAM029"

# Bring over the formatting used by the other "character length" sheets
# (IAM005/IAM006 share the exact same 4-column layout/styles).
$iam006 = $wb.Worksheets.Item("IAM006")
$iam006.Range("A1:D4").Copy()
$iam029.Range("A1:D4").PasteSpecial(-4122)

$iam029.Range("A1").Value = "CHARACTER LENGTH"
$iam029.Range("B1").Value = "VALIDITY"
$iam029.Range("C1").Value = "Runmode"
$iam029.Range("D1").Value = "PASS"

$iam029.Range("A2").Value = 91
$iam029.Range("B2").Value = "YES"
$iam029.Range("C2").Value = "Y"
$iam029.Range("D2").Value = "SKIP"

$iam029.Range("A3").Value = 92
$iam029.Range("B3").Value = "YES"
$iam029.Range("C3").Value = "Y"
$iam029.Range("D3").Value = "SKIP"

$iam029.Range("A4").Value = 93
$iam029.Range("B4").Value = "NO"
$iam029.Range("C4").Value = "Y"
$iam029.Range("D4").Value = "PASS"

$iam029.Range("A4").Select() | Out-Null

# ------------------------------------------------------------------
# 2. Add the new summary row (30) to the "Test Cases" sheet
# ------------------------------------------------------------------
$testCases = $wb.Worksheets.Item("Test Cases")

$testCases.Range("A29:E29").Copy()
$testCases.Range("A30:E30").PasteSpecial(-4122)

$testCases.Range("A30").Value = "IAM029"
$testCases.Range("B30").Value = "OPQA-2906"
$testCases.Range("C30").Value = "Verify that to validate PASSWORD field in new Neon user registration page with maximum length."
$testCases.Range("D30").Value = "Y"
$testCases.Range("E30").Value = "PASS"
$testCases.Rows.Item(30).RowHeight = 28.8

# ------------------------------------------------------------------
# 3. Update the selection on "IAM005" to the full A1:D4 range
# ------------------------------------------------------------------
$iam005 = $wb.Worksheets.Item("IAM005")
$iam005.Range("A1:D4").Select() | Out-Null

# ------------------------------------------------------------------
# 4. Restore "Test Cases" as the active sheet/cell (C30), matching the
#    unchanged tabSelected="1" on sheet1 in the diff.
# ------------------------------------------------------------------
$testCases.Select() | Out-Null
$testCases.Range("C30").Select() | Out-Null
